# ----------------------------------------------------------------------
# UP/Down Icons korrigiert: alle haben selbe Groesse.
#
# 1) Fix the fixed "date last changed" field text on the slide master
#    and on every slide layout (02.01.2015 -> 08.01.2015).
# 2) Rebuild the eight little "Pfeil nach rechts" (right arrow) shapes
#    that make up the up/down icon on slide 1 so that they all share
#    the exact same size and are laid out in two aligned columns.
# ----------------------------------------------------------------------

$p = $ppt.ActivePresentation

# Converts an EMU (English Metric Unit) length into points, the unit
# used by the Shape.Left/.Top properties. A tiny epsilon is added to
# counter-act float rounding that otherwise truncates the value by a
# single EMU.
function EmuToPt($emu) {
    return $emu / 12700.0 + 0.00002
}

# ------------------------------------------------------------------
# 1) Update the fixed date text "02.01.2015" -> "08.01.2015"
#    wherever the date placeholder shows up (slide master + all
#    custom layouts).
# ------------------------------------------------------------------
function Update-DateShape($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        if ($shp.Name -like "*Datum*") {
            if ($shp.TextFrame.TextRange.Text -eq "02.01.2015") {
                $shp.TextFrame.TextRange.Text = "08.01.2015"
            }
        }
    }
}

$design = $p.Designs.Item(1)
$master = $design.SlideMaster

Update-DateShape $master.Shapes

for ($li = 1; $li -le $master.CustomLayouts.Count; $li++) {
    $layout = $master.CustomLayouts.Item($li)
    Update-DateShape $layout.Shapes
}

# ------------------------------------------------------------------
# 2) Rebuild the 8 right-arrow shapes on slide 1.
# ------------------------------------------------------------------
$s = $p.Slides.Item(1)

# Use the first arrow as a template to duplicate from, since
# Duplicate() is the only way to create a new shape that keeps the
# original "p:style" (lnRef/fillRef/effectRef/fontRef) block and
# txBody structure intact.
$template = $s.Shapes.Item(1)

# Creating new shapes allocates new, ever increasing shape ids. Two
# throw-away duplicates are created and immediately deleted first so
# that the 8 real replacement shapes line up with the ids used by the
# target file (12..19).
$burn1 = $template.Duplicate().Item(1)
$burn1.Delete()
$burn2 = $template.Duplicate().Item(1)
$burn2.Delete()

# Target state for the 8 replacement shapes, in document order.
$targets = @(
    @{ name = "Pfeil nach rechts 11"; rot = 270; x = 2483768; y = 2348881; fill = "solid" },
    @{ name = "Pfeil nach rechts 12"; rot = 90;  x = 2292127; y = 2348881; fill = "solid" },
    @{ name = "Pfeil nach rechts 13"; rot = 270; x = 2483768; y = 2573288; fill = "none"  },
    @{ name = "Pfeil nach rechts 14"; rot = 90;  x = 2292127; y = 2573288; fill = "none"  },
    @{ name = "Pfeil nach rechts 15"; rot = 270; x = 2483768; y = 2786535; fill = "none"  },
    @{ name = "Pfeil nach rechts 16"; rot = 90;  x = 2292127; y = 2786535; fill = "solid" },
    @{ name = "Pfeil nach rechts 17"; rot = 270; x = 2483768; y = 3010942; fill = "solid" },
    @{ name = "Pfeil nach rechts 18"; rot = 90;  x = 2292127; y = 3010942; fill = "none"  }
)

foreach ($t in $targets) {
    $shp = $template.Duplicate().Item(1)
    $shp.Name = $t.name
    $shp.Rotation = $t.rot
    $shp.Left = EmuToPt $t.x
    $shp.Top = EmuToPt $t.y
    if ($t.fill -eq "solid") {
        $shp.Fill.Solid()
        $shp.Fill.ForeColor.ObjectThemeColor = 5
    } else {
        $shp.Fill.Visible = $false
    }
}

# Remove the original 8 arrow shapes (they are still the first 8
# shapes in the tree, in front of the freshly duplicated ones).
for ($i = 1; $i -le 8; $i++) {
    $s.Shapes.Item(1).Delete()
}
